$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.151.55'
$ws.Range("E2").Value = '  +1.69%  '

# Row 3
$ws.Range("D3").Value = '2.260.07'
$ws.Range("E3").Value = '  +3.17%  '

# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.44%  '

# Row 5
$ws.Range("D5").Value = '''249.16'
$ws.Range("E5").Value = '  -0.95%  '

# Row 6
$ws.Range("D6").Value = '''0.630'
$ws.Range("E6").Value = '  +2.27%  '

# Row 7
$ws.Range("D7").Value = '''70.86'
$ws.Range("E7").Value = '  +4.95%  '

# Row 8
$ws.Range("E8").Value = '  +0.24%  '

# Row 9
$ws.Range("D9").Value = '''0.638'
$ws.Range("E9").Value = '  +6.42%  '

# Row 10
$ws.Range("D10").Value = '''38.14'
$ws.Range("E10").Value = '  -3.78%  '

# Row 11
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = '''59.28'
$ws.Range("E11").Value = '  -0.80%  '

# Row 12
$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").Value = '''0.0959'
$ws.Range("E12").Value = '  +1.83%  '

# Row 13
$ws.Range("D13").Value = '''7.23'
$ws.Range("E13").Value = '  +3.47%  '

# Row 14
$ws.Range("D14").Value = '''0.105'
$ws.Range("E14").Value = '  +1.14%  '

# Row 15
$ws.Range("D15").Value = '2.594.50'
$ws.Range("E15").Value = '  +3.09%  '

# Row 16
$ws.Range("D16").Value = '''14.80'
$ws.Range("E16").Value = '  +2.09%  '

# Row 17
$ws.Range("D17").Value = '''0.872'
$ws.Range("E17").Value = '  +1.81%  '

# Row 18
$ws.Range("D18").Value = '2.278.26'
$ws.Range("E18").Value = '  +4.47%  '

# Row 19
$ws.Range("D19").Value = '42.224.73'
$ws.Range("E19").Value = '  +1.96%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0986'
$ws.Range("E20").Value = '  +3.67%  '

# Row 21
$ws.Range("E21").Value = '  +1.82%  '

# Row 22
$ws.Range("D22").Value = '''72.47'
$ws.Range("E22").Value = '  +0.67%  '

# Row 23
$ws.Range("D23").Value = '''2.26'
$ws.Range("E23").Value = '  +10.21%  '

# Row 24
$ws.Range("D24").Value = '''233.78'
$ws.Range("E24").Value = '  +1.27%  '

# Row 25
$ws.Range("D25").Value = '''3.91'
$ws.Range("E25").Value = '  +1.51%  '

# Row 26
$ws.Range("D26").Value = '''11.43'
$ws.Range("E26").Value = '  +0.38%  '

# Row 27
$ws.Range("D27").Value = '''0.997'
$ws.Range("E27").Value = '  -0.51%  '

# Row 28
$ws.Range("D28").Value = '''2.42'
$ws.Range("E28").Value = '  +0.20%  '

# Row 29
$ws.Range("D29").Value = '''3.63'
$ws.Range("E29").Value = '  -1.47%  '

# Row 30
$ws.Range("E30").Value = '  +2.12%  '

# Row 31
$ws.Range("D31").Value = '''166.29'
$ws.Range("E31").Value = '  -0.48%  '

# Row 32
$ws.Range("D32").Value = '''20.87'
$ws.Range("E32").Value = '  +2.97%  '

# Row 33
$ws.Range("D33").Value = '''6.31'
$ws.Range("E33").Value = '  +9.59%  '

# Row 34
$ws.Range("D34").Value = '''0.125'
$ws.Range("E34").Value = '  +3.95%  '

# Row 35
$ws.Range("E35").Value = '  +2.33%  '

# Row 36
$ws.Range("D36").Value = '''30.85'
$ws.Range("E36").Value = '  +19.36%  '

# Row 37
$ws.Range("E37").Value = '  +2.74%  '

# Row 38
$ws.Range("E38").Value = '  +10.24%  '

# Row 39
$ws.Range("E39").Value = '  +2.78%  '

# Row 40
$ws.Range("D40").Value = '''0.0305'
$ws.Range("E40").Value = '  +0.01%  '

# Row 41
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").Value = '''13.38'
$ws.Range("E41").Value = '  +13.94%  '

# Row 42
$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").Value = '''2.30'
$ws.Range("E42").Value = '  +3.55%  '

# Row 43
$ws.Range("D43").Value = '''5.90'
$ws.Range("E43").Value = '  +4.27%  '

# Row 44
$ws.Range("D44").Value = '''0.208'
$ws.Range("E44").Value = '  +7.77%  '

# Row 45
$ws.Range("D45").Value = '''9.17'
$ws.Range("E45").Value = '  +7.37%  '

# Row 46
$ws.Range("D46").Value = '''61.26'
$ws.Range("E46").Value = '  +0.08%  '

# Row 47
$ws.Range("D47").Value = '''4.86'
$ws.Range("E47").Value = '  -7.07%  '

# Row 48
$ws.Range("E48").Value = '  +3.13%  '

# Row 49
$ws.Range("E49").Value = '  +0.33%  '

# Row 50
$ws.Range("E50").Value = '  +0.27%  '

# Row 51
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value = '''1.17'
$ws.Range("E51").Value = '  +1.45%  '
